$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Satisfaction between Age Groups")

# --- New summary table data (rows 23-28) -------------------------------
$ws.Range("A23").Value = "Column1"
$ws.Range("B23").Value = "Average 18-30"
$ws.Range("C23").Value = "Average 50+"

$ws.Range("A24").Value = "Test1"
$ws.Range("B24").Formula = "=B20"
$ws.Range("C24").Formula = "=B9"

$ws.Range("A25").Value = "Test2"
$ws.Range("B25").Formula = "=C20"
$ws.Range("C25").Formula = "=C9"

$ws.Range("A26").Value = "Test3"
$ws.Range("B26").Formula = "=D20"
$ws.Range("C26").Formula = "=D9"

$ws.Range("A27").Value = "Test4"
$ws.Range("B27").Formula = "=E20"
$ws.Range("C27").Formula = "=E9"

$ws.Range("A28").Value = "Test5"
$ws.Range("B28").Formula = "=F20"
$ws.Range("C28").Formula = "=F9"

# Match the workbook's body font (Arial) used throughout this sheet so the
# new cells pick up the same style as the rest of the report.
$ws.Range("A23:F26").Font.Name = "Arial"
$ws.Range("A27:A28").Font.Name = "Arial"

# --- Turn the new range into a table (matches Table1..Table4 styling) --
$lo = $ws.ListObjects.Add(1, $ws.Range("A23:C28"), 0, 1)
$lo.Name = "Table6"
$lo.TableStyle = "TableStyleLight9"
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false

# --- Column widths roughly matching the new header text -----------------
$ws.Columns.Item(1).ColumnWidth = 11.5703125
$ws.Columns.Item(2).ColumnWidth = 16.5703125
$ws.Columns.Item(3).ColumnWidth = 14.85546875

# --- New line chart comparing the two averages across tests -------------
$chartObj = $ws.ChartObjects().Add(504824, 371475, 1543050, 2876550)
$chartObj.Name = "Chart 8"
$chart = $chartObj.Chart
$chart.ChartType = 4
$chart.SetSourceData($ws.Range("A23:C28"), 2)

$chart.HasTitle = $false

$series1 = $chart.SeriesCollection(1)
$series1.HasErrorBars = $true
$series2 = $chart.SeriesCollection(2)
$series2.HasErrorBars = $true

$xAxis = $chart.Axes(1)
$xAxis.HasTitle = $true
$xAxis.AxisTitle.Text = "Test Peformed"

$yAxis = $chart.Axes(2)
$yAxis.HasTitle = $true
$yAxis.AxisTitle.Text = "Average Satisfied from 0 to 10"
$yAxis.HasMajorGridlines = $true

$chart.HasLegend = $true
$chart.Legend.Position = -4152

# --- Selection / active cell matching the final editing state -----------
$ws.Activate()
$ws.Range("C27").Select()
